# Add a new paragraph with yellow-highlighted text right after the
# current final paragraph ("...In article on habr.") and before the
# section break, matching the author's "Added static library slides
# and example of throwing away unused code" commit.

$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$lastRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.Text = "Maybe find another name for hereditary disease"
$newRange.HighlightColorIndex = 7
